# Update "想去人数" (want-to-go count) figures for the refreshed data pull.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 10881
$wsExhibit.Range("F3").Value = 242
$wsExhibit.Range("F4").Value = 73
$wsExhibit.Range("F5").Value = 757
$wsExhibit.Range("F6").Value = 509

# Sheet "全部类型" (rId4 / sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10881
$wsAll.Range("F3").Value = 242
$wsAll.Range("F4").Value = 73
$wsAll.Range("F5").Value = 757
$wsAll.Range("F7").Value = 509
